$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Row 11 ("Marking"): Right count corrected from 5 to 4, Wrong marking from -1 to -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 ("Total"): recalc total marks and marks fraction display
$ws.Range("B12").Value = 92
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "90 / 112"
